$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Mark rows 41 through 70 (column B) as "ok" in the Status column
for ($r = 41; $r -le 70; $r++) {
    $ws.Cells.Item($r, 2).Value = "ok"
}

# Recalculate so dependent formulas (E2:E4) update
$excel.Calculate()

# Adjust the window view: split the pane (not frozen) and update selections
$ws.Activate()
$ws.Range("A61").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.Split = $true
$excel.ActiveWindow.SplitRow = 4
$excel.ActiveWindow.SplitColumn = 0

$ws.Range("D1:E4").Select()
$ws.Range("B70").Select()
